$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - PACOTE PRÉ-OPERATÓRIO PEDIÁTRICO OTORRINO
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1

# Row 3 - PACOTE PRÉ-OPERATÓRIO PEDIÁTRICO CIRURGIA GERAL
$ws.Range("H3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = 0

# Row 5 - ADENOIDECTOMIA PEDIÁTRICO
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1

# Row 6 - AMIGDALECTOMIA- PEDIATRICO
$ws.Range("D6").Value = 1

# Row 7 - AMIGDALECTOMIA COM ADENOIDECTOMIA - PEDIATRICO
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1

# Row 11 - HERNIOPLASTIA UMBILICAL - PEDIATRICO
$ws.Range("K11").Value = 0

# Row 16 - POSTECTOMIA - PEDIATRICO
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = 0

# Row 17 - TOTAL
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 3
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 2
$ws.Range("M17").Value = 0
